# Auto-generated Excel COM-interop script to update Sheets/Bahamut_Profits.xlsx (Leve profit tables)
# Applies cell value updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 772557.25
$ws.Range("J129").Value = 1059375.9
$ws.Range("L129").Value = 3178127.7
$ws.Range("N129").Value = -3188127.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 52601.5
$ws.Range("I2").Value = 79157.30499999999
$ws.Range("J2").Value = 3283.5715
$ws.Range("K2").Value = 79157.30499999999
$ws.Range("L2").Value = 3283.5715
$ws.Range("M2").Value = -79044.30499999999
$ws.Range("N2").Value = -3509.5715

$ws.Range("H39").Value = 12500
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21040

$ws.Range("H61").Value = 1480.1282
$ws.Range("I61").Value = 1518.7142
$ws.Range("J61").Value = 1142.5
$ws.Range("K61").Value = 1518.7142
$ws.Range("L61").Value = 1142.5
$ws.Range("M61").Value = -1306.7142
$ws.Range("N61").Value = -1566.5

$ws.Range("H74").Value = 1160.9474
$ws.Range("I74").Value = 1148.5454
$ws.Range("J74").Value = 1242.8
$ws.Range("K74").Value = 1148.5454
$ws.Range("L74").Value = 1242.8
$ws.Range("M74").Value = -274.5454
$ws.Range("N74").Value = -2990.8

$ws.Range("H77").Value = 1160.9474
$ws.Range("I77").Value = 1148.5454
$ws.Range("J77").Value = 1242.8
$ws.Range("K77").Value = 5742.727
$ws.Range("L77").Value = 6214
$ws.Range("M77").Value = -1374.727
$ws.Range("N77").Value = -14950

$ws.Range("H116").Value = 52601.5
$ws.Range("I116").Value = 79157.30499999999
$ws.Range("J116").Value = 3283.5715
$ws.Range("K116").Value = 79157.30499999999
$ws.Range("L116").Value = 3283.5715
$ws.Range("M116").Value = -76863.30499999999
$ws.Range("N116").Value = -7871.5715

$ws.Range("H122").Value = 955.1111
$ws.Range("I122").Value = 800.3333
$ws.Range("K122").Value = 2400.9999
$ws.Range("M122").Value = 49.0001000000002

$ws.Range("H132").Value = 1058.2153
$ws.Range("I132").Value = 915.9032
$ws.Range("K132").Value = 2747.7096
$ws.Range("M132").Value = -217.7096000000001

$ws.Range("H136").Value = 1480.1282
$ws.Range("I136").Value = 1518.7142
$ws.Range("J136").Value = 1142.5
$ws.Range("K136").Value = 4556.142599999999
$ws.Range("L136").Value = 3427.5
$ws.Range("M136").Value = -2006.142599999999
$ws.Range("N136").Value = -8527.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 52601.5
$ws.Range("I3").Value = 79157.30499999999
$ws.Range("J3").Value = 3283.5715
$ws.Range("K3").Value = 79157.30499999999
$ws.Range("L3").Value = 3283.5715
$ws.Range("M3").Value = -79043.30499999999
$ws.Range("N3").Value = -3511.5715

$ws.Range("H105").Value = 3701.9348
$ws.Range("I105").Value = 3461.3103
$ws.Range("J105").Value = 4112.4116
$ws.Range("K105").Value = 3461.3103
$ws.Range("L105").Value = 4112.4116
$ws.Range("M105").Value = -1714.3103
$ws.Range("N105").Value = -7606.4116

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 8000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 8000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1874.4
$ws.Range("I5").Value = 389.30768
$ws.Range("J5").Value = 3483.25
$ws.Range("K5").Value = 1167.92304
$ws.Range("L5").Value = 10449.75
$ws.Range("M5").Value = -1055.92304
$ws.Range("N5").Value = -10673.75

$ws.Range("H22").Value = 21717172
$ws.Range("I22").Value = 21717172
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 65151516
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -65151347
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 21717172
$ws.Range("I27").Value = 21717172
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 65151516
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -65151414
$ws.Range("N27").ClearContents()

$ws.Range("H98").Value = 2576.8
$ws.Range("I98").Value = 300
$ws.Range("J98").Value = 2829.7778
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 8489.3334
$ws.Range("M98").Value = 598
$ws.Range("N98").Value = -11485.3334

$ws.Range("H117").Value = 3024.077
$ws.Range("J117").Value = 3993.7778
$ws.Range("L117").Value = 11981.3334
$ws.Range("N117").Value = -18865.3334

$ws.Range("H121").Value = 926.6786
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()

$ws.Range("H135").Value = 1874.4
$ws.Range("I135").Value = 389.30768
$ws.Range("J135").Value = 3483.25
$ws.Range("K135").Value = 3503.76912
$ws.Range("L135").Value = 31349.25
$ws.Range("M135").Value = -968.7691199999999
$ws.Range("N135").Value = -36419.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 44538.46
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 44538.46
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 44538.46
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -46160.46

$ws.Range("H71").Value = 44538.46
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 44538.46
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 133615.38
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -141727.38

$ws.Range("H98").Value = 80000
$ws.Range("J98").Value = 80000
$ws.Range("L98").Value = 80000
$ws.Range("N98").Value = -85990

$ws.Range("H99").Value = 20000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 20000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -24492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1898.762
$ws.Range("I7").Value = 2056
$ws.Range("J7").Value = 1584.2858
$ws.Range("K7").Value = 2056
$ws.Range("L7").Value = 1584.2858
$ws.Range("M7").Value = -1944
$ws.Range("N7").Value = -1808.2858

$ws.Range("H93").Value = 3169.238
$ws.Range("I93").Value = 3027.7778
$ws.Range("J93").Value = 3275.3333
$ws.Range("K93").Value = 3027.7778
$ws.Range("L93").Value = 3275.3333
$ws.Range("M93").Value = -1779.7778
$ws.Range("N93").Value = -5771.3333

$ws.Range("H126").Value = 1898.762
$ws.Range("I126").Value = 2056
$ws.Range("J126").Value = 1584.2858
$ws.Range("K126").Value = 6168
$ws.Range("L126").Value = 4752.857400000001
$ws.Range("M126").Value = -3698
$ws.Range("N126").Value = -9692.857400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 528.9091
$ws.Range("I107").Value = 462.4643
$ws.Range("J107").Value = 901
$ws.Range("K107").Value = 1387.3929
$ws.Range("L107").Value = 2703
$ws.Range("M107").Value = 532.6071000000002
$ws.Range("N107").Value = -6543

Write-Host "Update complete"
